$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.634.04'
$ws.Range("E2").Value = '  -2.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.759.57'
$ws.Range("E3").Value = '  -3.32%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.53'
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4302'
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3605'
$ws.Range("E8").Value = '  -1.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07571'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.10'
$ws.Range("E10").Value = '  -6.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.110'
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.80'
$ws.Range("E13").Value = '  -6.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.071'
$ws.Range("E14").Value = '  -3.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.232'
$ws.Range("E15").Value = '  -4.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.760.66'
$ws.Range("E16").Value = '  -4.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.24'
$ws.Range("E17").Value = '  -1.71%  '
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06428'
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.07'
$ws.Range("E21").Value = '  -2.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.874'
$ws.Range("E22").Value = '  -6.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.694.71'
$ws.Range("E23").Value = '  -2.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.26'
$ws.Range("E24").Value = '  -3.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.097'
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.57'
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.959.13'
$ws.Range("E28").Value = '  -4.22%  '
$ws.Range("E29").Value = '  -6.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.87'
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.101'
$ws.Range("E31").Value = '  -10.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.685'
$ws.Range("E32").Value = '  +6.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.606'
$ws.Range("E33").Value = '  -6.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08958'
$ws.Range("E34").Value = '  -2.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.22'
$ws.Range("E35").Value = '  -6.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02305'
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2113'
$ws.Range("E37").Value = '  -3.23%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6395'
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06016'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.950'
$ws.Range("E40").Value = '  -4.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.189'
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.397'
$ws.Range("E43").Value = '  -2.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.930'
$ws.Range("E44").Value = '  -2.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.41'
$ws.Range("E45").Value = '  -3.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5928'
$ws.Range("E46").Value = '  -3.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.712'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.988'
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.57'
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.168'
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("E51").Value = '  -1.93%  '
